$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "Data Science 대학원 수학&통계학 Boot Camp 난이도"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/boot-camp-difficulty/#utm_source=rss&utm_medium=rss&utm_campaign=boot-camp-difficulty"

$ws.Range("D28").Value = "[임피던스 제어] 상호작용 컨트롤(1)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/101"

$ws.Range("D32").Value = "Edit Distance (Levenshtein Distance) (퍼옴)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/306"

$ws.Range("D36").Value = "Introduction to Image Super Resolution"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/318"

$ws.Range("D37").Value = "[paper Review] Contrastive Learning for Sequential Recommendation"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1454&mod=document&pageid=1"

$ws.Range("D39").Value = "Using tf.Print() in TensorFlow"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Using-tfPrint-in-TensorFlow-1"

$ws.Range("D50").Value = "컴파일러의 탄생"
$ws.Range("E50").Value = "http://incredible.egloos.com/7514515"

$ws.Range("D51").Value = "[python] 파일의 확장자를 알려주는 함수, os.path.splitext()"
$ws.Range("E51").Value = "https://bskyvision.com/1077"

$wb.Save()
